# Append three more config rows (20-22) for US/UK/JP with a Tracking
# Error Constraint of "Yes" and a 5% Tracking Error Limit, mirroring the
# first three country rows but with the new constraint columns filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRows = @(
    @{ Row = 21; Config = 20; Country = "US" },
    @{ Row = 22; Config = 21; Country = "UK" },
    @{ Row = 23; Config = 22; Country = "JP" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Config      # A: Config
    $ws.Cells.Item($row, 2).Value = $r.Country     # B: Country
    $ws.Cells.Item($row, 3).Value = 0.07            # C: Target Volatility
    $ws.Cells.Item($row, 4).Value = 0.1             # D: Epsilon
    $ws.Cells.Item($row, 5).Value = 40              # E: Lambda
    $ws.Cells.Item($row, 6).Value = "None"          # F: Additional Constraints
    $ws.Cells.Item($row, 7).Value = "Yes"           # G: Tracking Error Constraint
    $ws.Cells.Item($row, 8).Value = 0.05             # H: Tracking Error Limit
}

# Mirror the final cursor/selection state recorded in the workbook.
$ws.Range("X26").Select() | Out-Null
